$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.950.42'
$ws.Range('E2').Value = '  -3.89%  '
$ws.Range('D3').Value = '3.510.56'
$ws.Range('E3').Value = '  -4.85%  '
$ws.Range('D4').NumberFormatLocal = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').NumberFormatLocal = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormatLocal = '@'
$ws.Range('D5').Value = '579.66'
$ws.Range('D5').NumberFormatLocal = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').NumberFormatLocal = '@'
$ws.Range('D6').Value = '174.62'
$ws.Range('D6').NumberFormatLocal = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.501.68'
$ws.Range('E8').Value = '  -4.92%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormatLocal = '@'
$ws.Range('D10').Value = '0.189'
$ws.Range('D10').NumberFormatLocal = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.78%  '
$ws.Range('D11').NumberFormatLocal = '@'
$ws.Range('D11').Value = '6.70'
$ws.Range('D11').NumberFormatLocal = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +9.50%  '
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D13').NumberFormatLocal = '@'
$ws.Range('D13').Value = '47.28'
$ws.Range('D13').NumberFormatLocal = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.52%  '
$ws.Range('D14').NumberFormatLocal = '@'
$ws.Range('D14').Value = '0.0000277'
$ws.Range('D14').NumberFormatLocal = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.27%  '
$ws.Range('D15').NumberFormatLocal = '@'
$ws.Range('D15').Value = '671.96'
$ws.Range('D15').NumberFormatLocal = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').Value = '4.063.74'
$ws.Range('E16').Value = '  -5.22%  '
$ws.Range('D17').NumberFormatLocal = '@'
$ws.Range('D17').Value = '8.83'
$ws.Range('D17').NumberFormatLocal = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '3.500.09'
$ws.Range('E18').Value = '  -4.95%  '
$ws.Range('D19').Value = '68.821.43'
$ws.Range('E19').Value = '  -4.34%  '
$ws.Range('E20').Value = '  -1.60%  '
$ws.Range('D21').NumberFormatLocal = '@'
$ws.Range('D21').Value = '17.57'
$ws.Range('D21').NumberFormatLocal = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.48%  '
$ws.Range('D22').NumberFormatLocal = '@'
$ws.Range('D22').Value = '11.20'
$ws.Range('D22').NumberFormatLocal = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.74%  '
$ws.Range('D23').NumberFormatLocal = '@'
$ws.Range('D23').Value = '0.906'
$ws.Range('D23').NumberFormatLocal = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.74%  '
$ws.Range('D24').NumberFormatLocal = '@'
$ws.Range('D24').Value = '16.29'
$ws.Range('D24').NumberFormatLocal = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.29%  '
$ws.Range('D25').NumberFormatLocal = '@'
$ws.Range('D25').Value = '98.34'
$ws.Range('D25').NumberFormatLocal = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.00%  '
$ws.Range('D26').NumberFormatLocal = '@'
$ws.Range('D26').Value = '3.87'
$ws.Range('D26').NumberFormatLocal = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.23%  '
$ws.Range('D27').NumberFormatLocal = '@'
$ws.Range('D27').Value = '5.82'
$ws.Range('D27').NumberFormatLocal = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('E28').Value = '  +0.06%  '
$ws.Range('D29').NumberFormatLocal = '@'
$ws.Range('D29').Value = '2.65'
$ws.Range('D29').NumberFormatLocal = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.03%  '
$ws.Range('D30').NumberFormatLocal = '@'
$ws.Range('D30').Value = '9.44'
$ws.Range('D30').NumberFormatLocal = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.47%  '
$ws.Range('D31').NumberFormatLocal = '@'
$ws.Range('D31').Value = '32.96'
$ws.Range('D31').NumberFormatLocal = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.16%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormatLocal = '@'
$ws.Range('D32').Value = '3.22'
$ws.Range('D32').NumberFormatLocal = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.50%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormatLocal = '@'
$ws.Range('D33').Value = '8.74'
$ws.Range('D33').NumberFormatLocal = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.93%  '
$ws.Range('D34').NumberFormatLocal = '@'
$ws.Range('D34').Value = '7.31'
$ws.Range('D34').NumberFormatLocal = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').NumberFormatLocal = '@'
$ws.Range('D35').Value = '1.36'
$ws.Range('D35').NumberFormatLocal = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.64%  '
$ws.Range('D36').NumberFormatLocal = '@'
$ws.Range('D36').Value = '577.60'
$ws.Range('D36').NumberFormatLocal = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').NumberFormatLocal = '@'
$ws.Range('D37').Value = '10.94'
$ws.Range('D37').NumberFormatLocal = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.35%  '
$ws.Range('D38').NumberFormatLocal = '@'
$ws.Range('D38').Value = '3.58'
$ws.Range('D38').NumberFormatLocal = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -14.50%  '
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('D40').NumberFormatLocal = '@'
$ws.Range('D40').Value = '56.95'
$ws.Range('D40').NumberFormatLocal = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.85%  '
$ws.Range('D41').NumberFormatLocal = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').NumberFormatLocal = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').NumberFormatLocal = '@'
$ws.Range('D42').Value = '0.137'
$ws.Range('D42').NumberFormatLocal = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.33%  '
$ws.Range('D43').NumberFormatLocal = '@'
$ws.Range('D43').Value = '0.337'
$ws.Range('D43').NumberFormatLocal = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.12%  '
$ws.Range('D44').NumberFormatLocal = '@'
$ws.Range('D44').Value = '0.0438'
$ws.Range('D44').NumberFormatLocal = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.50%  '
$ws.Range('D45').Value = '3.421.35'
$ws.Range('E45').Value = '  -8.81%  '
$ws.Range('D46').NumberFormatLocal = '@'
$ws.Range('D46').Value = '33.44'
$ws.Range('D46').NumberFormatLocal = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.80%  '
$ws.Range('D47').Value = '0.0₃0706'
$ws.Range('E47').Value = '  -8.59%  '
$ws.Range('D48').NumberFormatLocal = '@'
$ws.Range('D48').Value = '2.91'
$ws.Range('D48').NumberFormatLocal = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('E49').Value = '  -7.09%  '
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').NumberFormatLocal = '@'
$ws.Range('D51').Value = '130.48'
$ws.Range('D51').NumberFormatLocal = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.62%  '
